# Weekly update: insert a new data row at row 133 (pushing existing rows
# 133-156 down to 134-157), and populate the newly inserted row with the
# latest week's price observation for Perejil at Feria Lagunitas de Puerto
# Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 133; this shifts rows 133:156
# down to 134:157, carrying all of their existing values/formatting along.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly record.
$ws.Range("A133").Value = 4
$ws.Range("B133").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C133").Value = "Los Lagos"
$ws.Range("D133").Value2 = 44476
$ws.Range("E133").Value = 10
$ws.Range("F133").Value = 100112044
$ws.Range("G133").Value = "Perejil"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 90
$ws.Range("K133").Value = 4500
$ws.Range("L133").Value = 4500
$ws.Range("M133").Value = 4500
$ws.Range("N133").Value = "$/docena de atados (3 kilos)"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 1500
$ws.Range("Q133").Value = 3
$ws.Range("R133").Value = "Hortaliza"
